$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before row 569, pushing the existing rows
# 569-574 (the 2021-11-15 / 2021-11-08 "Packham's Triumph" data) down to
# become rows 574-579, unchanged.
$ws.Rows("569:573").Insert(-4121)  # xlShiftDown

# Copy formatting (date number format on column D) from the row that is
# now 574 into the 5 newly inserted rows so D569:D573 render as dates too.
$ws.Range("D569:D573").NumberFormat = $ws.Range("D574").NumberFormat

# New week (2022-04-05, serial 44656) of "Pera" price data for
# "Terminal La Palmera de La Serena".

# Row 569: Packham's Triumph - Especial
$ws.Range("A569").Value = 8
$ws.Range("B569").Value = "Terminal La Palmera de La Serena"
$ws.Range("C569").Value = "Coquimbo"
$ws.Range("D569").Value = 44656
$ws.Range("E569").Value = 4
$ws.Range("F569").Value = "Fruta"
$ws.Range("G569").Value = 100104
$ws.Range("H569").Value = "Frutos de pepita"
$ws.Range("I569").Value = 100104005
$ws.Range("J569").Value = "Pera"
$ws.Range("K569").Value = "Packham's Triumph"
$ws.Range("L569").Value = "Especial"
$ws.Range("M569").Value = 20
$ws.Range("N569").Value = 245000
$ws.Range("O569").Value = 250000
$ws.Range("P569").Value = 247500
$ws.Range("Q569").Value = "`$/bins (450 kilos)"
$ws.Range("R569").Value = "Región de O'Higgins"
$ws.Range("S569").Value = 550
$ws.Range("T569").Value = 450

# Row 570: Packham's Triumph - Primera
$ws.Range("A570").Value = 8
$ws.Range("B570").Value = "Terminal La Palmera de La Serena"
$ws.Range("C570").Value = "Coquimbo"
$ws.Range("D570").Value = 44656
$ws.Range("E570").Value = 4
$ws.Range("F570").Value = "Fruta"
$ws.Range("G570").Value = 100104
$ws.Range("H570").Value = "Frutos de pepita"
$ws.Range("I570").Value = 100104005
$ws.Range("J570").Value = "Pera"
$ws.Range("K570").Value = "Packham's Triumph"
$ws.Range("L570").Value = "Primera"
$ws.Range("M570").Value = 16
$ws.Range("N570").Value = 215000
$ws.Range("O570").Value = 220000
$ws.Range("P570").Value = 217500
$ws.Range("Q570").Value = "`$/bins (450 kilos)"
$ws.Range("R570").Value = "Región de O'Higgins"
$ws.Range("S570").Value = 483
$ws.Range("T570").Value = 450

# Row 571: Packham's Triumph - Segunda
$ws.Range("A571").Value = 8
$ws.Range("B571").Value = "Terminal La Palmera de La Serena"
$ws.Range("C571").Value = "Coquimbo"
$ws.Range("D571").Value = 44656
$ws.Range("E571").Value = 4
$ws.Range("F571").Value = "Fruta"
$ws.Range("G571").Value = 100104
$ws.Range("H571").Value = "Frutos de pepita"
$ws.Range("I571").Value = 100104005
$ws.Range("J571").Value = "Pera"
$ws.Range("K571").Value = "Packham's Triumph"
$ws.Range("L571").Value = "Segunda"
$ws.Range("M571").Value = 14
$ws.Range("N571").Value = 195000
$ws.Range("O571").Value = 200000
$ws.Range("P571").Value = 197500
$ws.Range("Q571").Value = "`$/bins (450 kilos)"
$ws.Range("R571").Value = "Región de O'Higgins"
$ws.Range("S571").Value = 439
$ws.Range("T571").Value = 450

# Row 572: Winter Nelis - Especial
$ws.Range("A572").Value = 8
$ws.Range("B572").Value = "Terminal La Palmera de La Serena"
$ws.Range("C572").Value = "Coquimbo"
$ws.Range("D572").Value = 44656
$ws.Range("E572").Value = 4
$ws.Range("F572").Value = "Fruta"
$ws.Range("G572").Value = 100104
$ws.Range("H572").Value = "Frutos de pepita"
$ws.Range("I572").Value = 100104005
$ws.Range("J572").Value = "Pera"
$ws.Range("K572").Value = "Winter Nelis"
$ws.Range("L572").Value = "Especial"
$ws.Range("M572").Value = 20
$ws.Range("N572").Value = 275000
$ws.Range("O572").Value = 280000
$ws.Range("P572").Value = 277500
$ws.Range("Q572").Value = "`$/bins (450 kilos)"
$ws.Range("R572").Value = "Región de O'Higgins"
$ws.Range("S572").Value = 617
$ws.Range("T572").Value = 450

# Row 573: Winter Nelis - Primera
$ws.Range("A573").Value = 8
$ws.Range("B573").Value = "Terminal La Palmera de La Serena"
$ws.Range("C573").Value = "Coquimbo"
$ws.Range("D573").Value = 44656
$ws.Range("E573").Value = 4
$ws.Range("F573").Value = "Fruta"
$ws.Range("G573").Value = 100104
$ws.Range("H573").Value = "Frutos de pepita"
$ws.Range("I573").Value = 100104005
$ws.Range("J573").Value = "Pera"
$ws.Range("K573").Value = "Winter Nelis"
$ws.Range("L573").Value = "Primera"
$ws.Range("M573").Value = 16
$ws.Range("N573").Value = 225000
$ws.Range("O573").Value = 230000
$ws.Range("P573").Value = 227500
$ws.Range("Q573").Value = "`$/bins (450 kilos)"
$ws.Range("R573").Value = "Región de O'Higgins"
$ws.Range("S573").Value = 506
$ws.Range("T573").Value = 450
